$wb = $excel.ActiveWorkbook

# --- Sheet "gens" (sheet1): production column (C) updates ---
$wsGens = $wb.Worksheets.Item("gens")
$wsGens.Range("C10").Value = 82.5
$wsGens.Range("C11").Value = 100
$wsGens.Range("C12").Value = 100
$wsGens.Range("C22").Value = 155
$wsGens.Range("C23").Value = 155
$wsGens.Range("C25").Value = 400

# --- Sheet "lines" (sheet2): flow_loads/susceptance_shadow/flow_low_shadow/flow_high_shadow updates ---
$wsLines = $wb.Worksheets.Item("lines")
$wsLines.Range("C2").Value = 24.990052
$wsLines.Range("D2").Value = 1.9332613
$wsLines.Range("C3").Value = -91.732021
$wsLines.Range("D3").Value = -7.1995755
$wsLines.Range("C4").Value = -41.25803
$wsLines.Range("D4").Value = -8.837354700000001
$wsLines.Range("C5").Value = -50.509948
$wsLines.Range("D5").Value = -7.469632
$wsLines.Range("C6").Value = -21.5
$wsLines.Range("D6").Value = 37.805969
$wsLines.Range("C7").Value = 2.8699097
$wsLines.Range("D7").Value = -4.4536593
$wsLines.Range("C8").Value = -98.74149199999999
$wsLines.Range("D8").Value = 0.27757779
$wsLines.Range("C9").Value = -124.50995
$wsLines.Range("D9").Value = -6.116864
$wsLines.Range("C10").Value = -112.25803
$wsLines.Range("D10").Value = -9.149261299999999
$wsLines.Range("C11").Value = -157.5
$wsLines.Range("D11").Value = 12.011271
$wsLines.Range("E11").Value = 69.737118
$wsLines.Range("C12").Value = 157.5
$wsLines.Range("F12").Value = -1.6723225
$wsLines.Range("C13").Value = -23.805448
$wsLines.Range("D13").Value = 3.1666906
$wsLines.Range("C14").Value = 10.305448
$wsLines.Range("D14").Value = -3.1666906
$wsLines.Range("C15").Value = -112.75279
$wsLines.Range("D15").Value = -2.0826203
$wsLines.Range("C16").Value = -207.69269
$wsLines.Range("D16").Value = -4.3895498
$wsLines.Range("C17").Value = -179.75634
$wsLines.Range("D17").Value = 4.2507609
$wsLines.Range("C18").Value = -274.69624
$wsLines.Range("D18").Value = 1.9438314
$wsLines.Range("C19").Value = -356.42422
$wsLines.Range("D19").Value = -2.1027283
$wsLines.Range("C20").Value = 63.915082
$wsLines.Range("D20").Value = 2.9239575
$wsLines.Range("C21").Value = -190.27939
$wsLines.Range("D21").Value = 0.20420122
$wsLines.Range("C22").Value = -292.10954
$wsLines.Range("D22").Value = -3.2368791
$wsLines.Range("C23").Value = -220.70361
$wsLines.Range("D23").Value = -3.4410803
$wsLines.Range("C24").Value = -130.08492
$wsLines.Range("D24").Value = 4.1074642
$wsLines.Range("C25").Value = -40.673483
$wsLines.Range("D25").Value = 0.045730447
$wsLines.Range("C26").Value = -66.94219699999999
$wsLines.Range("D26").Value = 0.015054545
$wsLines.Range("C27").Value = -66.94219699999999
$wsLines.Range("D27").Value = 0.015054545
$wsLines.Range("C28").Value = 98.74149199999999
$wsLines.Range("D28").Value = -0.17183387
$wsLines.Range("C29").Value = 334.2416
$wsLines.Range("D29").Value = -0.015976252
$wsLines.Range("C30").Value = -450
$wsLines.Range("D30").Value = 1.6772183
$wsLines.Range("E30").Value = 16.513047
$wsLines.Range("C31").Value = -450
$wsLines.Range("D31").Value = -0.00741303
$wsLines.Range("E31").Value = 0.00040309249
$wsLines.Range("C32").Value = -183.11561
$wsLines.Range("D32").Value = -0.008921753100000001
$wsLines.Range("C33").Value = -191.5
$wsLines.Range("D33").Value = -0.0068835278
$wsLines.Range("C34").Value = -191.5
$wsLines.Range("D34").Value = -0.0068835278
$wsLines.Range("C35").Value = -315.5
$wsLines.Range("D35").Value = 1.4584507
$wsLines.Range("C36").Value = -315.5
$wsLines.Range("D36").Value = 1.4584507
$wsLines.Range("C37").Value = -379.5
$wsLines.Range("D37").Value = 0.80214788
$wsLines.Range("C38").Value = -379.5
$wsLines.Range("D38").Value = 0.80214788
$wsLines.Range("C39").Value = -116.88439
$wsLines.Range("D39").Value = 0.0057778972

# --- Sheet "bus" (sheet3): bus_lmp/node_theta updates ---
$wsBus = $wb.Worksheets.Item("bus")
$wsBus.Range("B2").Value = 116.49225
$wsBus.Range("C2").Value = -0.10227434
$wsBus.Range("B3").Value = 118.42551
$wsBus.Range("C3").Value = -0.10577295
$wsBus.Range("B4").Value = 109.29267
$wsBus.Range("C4").Value = 0.09128022700000001
$wsBus.Range("B5").Value = 110.95588
$wsBus.Range("C5").Value = -0.041625312
$wsBus.Range("B6").Value = 107.65489
$wsBus.Range("C6").Value = -0.06720501299999999
$wsBus.Range("B7").Value = 156.23148
$wsBus.Range("C7").Value = -0.064492946
$wsBus.Range("C8").Value = 0.14466104
$wsBus.Range("B9").Value = 101.67232
$wsBus.Range("C9").Value = 0.048586044
$wsBus.Range("B10").Value = 104.83901
$wsBus.Range("C10").Value = 0.08786503399999999
$wsBus.Range("B11").Value = 98.50563200000001
$wsBus.Range("C11").Value = 0.031582054
$wsBus.Range("B12").Value = 102.75639
$wsBus.Range("C12").Value = 0.18257738
$wsBus.Range("B13").Value = 100.44946
$wsBus.Range("C13").Value = 0.2623269
$wsBus.Range("B14").Value = 100.65366
$wsBus.Range("C14").Value = 0.35366101
$wsBus.Range("B15").Value = 105.68035
$wsBus.Range("C15").Value = 0.15573305
$wsBus.Range("B16").Value = 109.74208
$wsBus.Range("C16").Value = 0.22556866
$wsBus.Range("B17").Value = 109.78781
$wsBus.Range("C17").Value = 0.23248315
$wsBus.Range("B18").Value = 109.77184
$wsBus.Range("C18").Value = 0.14558033
$wsBus.Range("B19").Value = 109.76402
$wsBus.Range("C19").Value = 0.20858033
$wsBus.Range("B20").Value = 94.95198600000001
$wsBus.Range("C20").Value = 0.33598315
$wsBus.Range("B21").Value = 96.410436
$wsBus.Range("C21").Value = 0.46218315
$wsBus.Range("B22").Value = 109.75714
$wsBus.Range("C22").Value = 0.25837033
$wsBus.Range("B23").Value = 109.76292
$wsBus.Range("C23").Value = 0.33785172
$wsBus.Range("B24").Value = 97.21258400000001
$wsBus.Range("C24").Value = 0.5456731500000001
$wsBus.Range("B25").Value = 109.57025
$wsBus.Range("C25").Value = 0.17422308

Write-Host "Applied all updates"
